# Repull data, push all data, mean calculation
# Updates column F (dSF) values for rows 2-35 (row 26 / 36 unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 1
    4  = 4
    5  = -1
    6  = 2
    7  = 2
    8  = -2
    9  = -1
    10 = -1
    11 = -1
    12 = -1
    13 = -3
    14 = -2
    15 = 1
    16 = -2
    17 = -2
    18 = -3
    19 = 3
    20 = -4
    21 = -2
    22 = 4
    23 = 2
    24 = -1
    25 = -1
    27 = -3
    28 = -4
    29 = 5
    30 = 3
    31 = -4
    32 = 1
    33 = 1
    34 = -1
    35 = -3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
